# Fixed update to excel issue
#
# 1. Rename the "Requested quantity" header on the two existing sheets.
# 2. Add a new "PO Forecast" sheet (after the existing sheets) with a
#    Prophet-style forecast table: ds / PO_Forecast / yhat_lower / yhat_upper.

$wb = $excel.ActiveWorkbook

# --- 1. Rename headers on the existing sheets -----------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the "PO Forecast" sheet at the end -----------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Copy the header/date formatting (bold+border style, date number format)
# from the "Weekly Quantity" sheet so the new sheet reuses the same styles
# rather than creating new ones.
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A16").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Headers
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Data rows
$rows = @(
    @(45144.99999999999, 25, -30.92558061550366, 76.49378777132839),
    @(45158.99999999999, 27, -26.08934420527964, 83.53243886337668),
    @(45403.99999999999, 60, 7.397927265350193, 109.7899067501158),
    @(45410.99999999999, 60, 3.296946544343354, 117.8248846802474),
    @(45445.99999999999, 65, 9.849218670139605, 114.9928625691831),
    @(45459.99999999999, 67, 12.42001171072675, 121.3246374254171),
    @(45480.99999999999, 70, 11.70850059367522, 125.0567441906055),
    @(45487.99999999999, 71, 13.9988257384304, 124.4062834702497),
    @(45494.99999999999, 72, 18.9841960607223, 123.0914762408838),
    @(45501.99999999999, 73, 18.25310644493962, 122.3919042730715),
    @(45508.99999999999, 74, 21.69015024409917, 127.8574588004708),
    @(45515.99999999999, 74, 27.69594080909949, 130.094687810353),
    @(45522.99999999999, 75, 24.70156614760692, 132.227542252323),
    @(45529.99999999999, 76, 25.67979765296079, 131.8407173014228),
    @(45536.99999999999, 77, 22.3412137731565, 124.6198980300919)
)

$r = 2
foreach ($row in $rows) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
